# Updated run for publication
# Refresh the frequency-table proportions (rows 2-5, columns B-X) on Sheet1
# with the newly computed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00254452926208651
$ws.Range("C2").Value = 0.00127226463104326
$ws.Range("D2").Value = 0.964376590330789
$ws.Range("E2").Value = 0.00381679389312977
$ws.Range("F2").Value = 0.0190839694656489
$ws.Range("G2").Value = 0.0279898218829517
$ws.Range("H2").Value = 0.624681933842239
$ws.Range("I2").Value = 0.00127226463104326
$ws.Range("J2").Value = 0.0190839694656489
$ws.Range("K2").Value = 0.991094147582697
$ws.Range("L2").Value = 0.00508905852417303
$ws.Range("M2").Value = 0.982188295165394
$ws.Range("O2").Value = 0.786259541984733
$ws.Range("Q2").Value = 0.00254452926208651
$ws.Range("R2").Value = 0.0178117048346056
$ws.Range("S2").Value = 0.860050890585242
$ws.Range("T2").Value = 0.0127226463104326
$ws.Range("U2").Value = 0.00127226463104326
$ws.Range("V2").Value = 0.0381679389312977
$ws.Range("W2").Value = 0.834605597964377
$ws.Range("X2").Value = 0.00636132315521628
$ws.Range("B3").Value = 0.977099236641221
$ws.Range("C3").Value = 0.974554707379135
$ws.Range("D3").Value = 0.00381679389312977
$ws.Range("E3").Value = 0.00381679389312977
$ws.Range("F3").Value = 0.00254452926208651
$ws.Range("G3").Value = 0.00254452926208651
$ws.Range("H3").Value = 0.0674300254452926
$ws.Range("I3").Value = 0.371501272264631
$ws.Range("J3").Value = 0.909669211195929
$ws.Range("K3").Value = 0.00508905852417303
$ws.Range("L3").Value = 0.00127226463104326
$ws.Range("N3").Value = 0.00127226463104326
$ws.Range("O3").Value = 0.0190839694656489
$ws.Range("P3").Value = 0.187022900763359
$ws.Range("Q3").Value = 0.972010178117048
$ws.Range("R3").Value = 0.975826972010178
$ws.Range("S3").Value = 0.110687022900763
$ws.Range("T3").Value = 0.924936386768448
$ws.Range("U3").Value = 0.0229007633587786
$ws.Range("W3").Value = 0.00254452926208651
$ws.Range("X3").Value = 0.837150127226463
$ws.Range("B4").Value = 0.0152671755725191
$ws.Range("C4").Value = 0.0229007633587786
$ws.Range("D4").Value = 0.00636132315521628
$ws.Range("E4").Value = 0.970737913486005
$ws.Range("F4").Value = 0.977099236641221
$ws.Range("G4").Value = 0.968193384223919
$ws.Range("H4").Value = 0.306615776081425
$ws.Range("K4").Value = 0.00254452926208651
$ws.Range("L4").Value = 0.993638676844784
$ws.Range("M4").Value = 0.0178117048346056
$ws.Range("N4").Value = 0.00254452926208651
$ws.Range("O4").Value = 0.0190839694656489
$ws.Range("P4").Value = 0.0139949109414758
$ws.Range("Q4").Value = 0.0203562340966921
$ws.Range("R4").Value = 0.00508905852417303
$ws.Range("S4").Value = 0.0241730279898219
$ws.Range("T4").Value = 0.00127226463104326
$ws.Range("U4").Value = 0.00127226463104326
$ws.Range("V4").Value = 0.932569974554707
$ws.Range("W4").Value = 0.0241730279898219
$ws.Range("X4").Value = 0.146310432569975
$ws.Range("B5").Value = 0.00508905852417303
$ws.Range("C5").Value = 0.00127226463104326
$ws.Range("D5").Value = 0.0254452926208651
$ws.Range("E5").Value = 0.0216284987277354
$ws.Range("H5").Value = 0.00127226463104326
$ws.Range("I5").Value = 0.627226463104326
$ws.Range("J5").Value = 0.0687022900763359
$ws.Range("N5").Value = 0.99618320610687
$ws.Range("O5").Value = 0.175572519083969
$ws.Range("P5").Value = 0.798982188295165
$ws.Range("Q5").Value = 0.00381679389312977
$ws.Range("R5").Value = 0.00127226463104326
$ws.Range("S5").Value = 0.00508905852417303
$ws.Range("T5").Value = 0.0610687022900763
$ws.Range("U5").Value = 0.974554707379135
$ws.Range("V5").Value = 0.0292620865139949
$ws.Range("W5").Value = 0.138676844783715
$ws.Range("X5").Value = 0.0101781170483461
